# Generate Report for Handback
# Marks the two localized files as handed back (in sync with en-US) for
# both the zh-cn and de-de language sheets: fills in the "Latest Target
# File" / "Latest Handback File" columns (F/G) with hyperlinks, stamps a
# handback datetime, and flips the Status column to the handed-back text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet - summary status columns for both languages
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column (shared by both data rows)
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

# Row 2 (29be957d-...)
$wsZh.Range("F2").Value = "29be957d-ffa9-4113-ba39-e45210c6289d.md"
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/d48e5211f45ca9d03dbdbc84089abf20aca7b4b6/e2e/29be957d-ffa9-4113-ba39-e45210c6289d.md", "", "", "29be957d-ffa9-4113-ba39-e45210c6289d.md") | Out-Null

$wsZh.Range("G2").Value = "29be957d-ffa9-4113-ba39-e45210c6289d.fe69d0faef49cc3c2d473d0b89da99b00d1183db.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/30ebab788a676a64bcdaef64de39967492141a03/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/29be957d-ffa9-4113-ba39-e45210c6289d.fe69d0faef49cc3c2d473d0b89da99b00d1183db.zh-cn.xlf", "", "", "29be957d-ffa9-4113-ba39-e45210c6289d.fe69d0faef49cc3c2d473d0b89da99b00d1183db.zh-cn.xlf") | Out-Null

$wsZh.Range("H2").Value = "2016-03-19 00:47:35"

# Row 3 (3cc8a830-...)
$wsZh.Range("F3").Value = "3cc8a830-1acc-4cdd-97cf-07c84515db63.md"
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/d48e5211f45ca9d03dbdbc84089abf20aca7b4b6/e2e/3cc8a830-1acc-4cdd-97cf-07c84515db63.md", "", "", "3cc8a830-1acc-4cdd-97cf-07c84515db63.md") | Out-Null

$wsZh.Range("G3").Value = "3cc8a830-1acc-4cdd-97cf-07c84515db63.960381a50bbb9142e9ab989eaff9b9e825477b0b.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/30ebab788a676a64bcdaef64de39967492141a03/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3cc8a830-1acc-4cdd-97cf-07c84515db63.960381a50bbb9142e9ab989eaff9b9e825477b0b.zh-cn.xlf", "", "", "3cc8a830-1acc-4cdd-97cf-07c84515db63.960381a50bbb9142e9ab989eaff9b9e825477b0b.zh-cn.xlf") | Out-Null

$wsZh.Range("H3").Value = "2016-03-19 00:47:35"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status column (shared by both data rows)
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# Row 2 (29be957d-...)
$wsDe.Range("F2").Value = "29be957d-ffa9-4113-ba39-e45210c6289d.md"
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/d48e5211f45ca9d03dbdbc84089abf20aca7b4b6/e2e/29be957d-ffa9-4113-ba39-e45210c6289d.md", "", "", "29be957d-ffa9-4113-ba39-e45210c6289d.md") | Out-Null

$wsDe.Range("G2").Value = "29be957d-ffa9-4113-ba39-e45210c6289d.fe69d0faef49cc3c2d473d0b89da99b00d1183db.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b71f158a6df8bc794c22d9aa2ebe3f2247e0ab85/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/29be957d-ffa9-4113-ba39-e45210c6289d.fe69d0faef49cc3c2d473d0b89da99b00d1183db.de-de.xlf", "", "", "29be957d-ffa9-4113-ba39-e45210c6289d.fe69d0faef49cc3c2d473d0b89da99b00d1183db.de-de.xlf") | Out-Null

$wsDe.Range("H2").Value = "2016-03-19 00:47:40"

# Row 3 (3cc8a830-...)
$wsDe.Range("F3").Value = "3cc8a830-1acc-4cdd-97cf-07c84515db63.md"
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/d48e5211f45ca9d03dbdbc84089abf20aca7b4b6/e2e/3cc8a830-1acc-4cdd-97cf-07c84515db63.md", "", "", "3cc8a830-1acc-4cdd-97cf-07c84515db63.md") | Out-Null

$wsDe.Range("G3").Value = "3cc8a830-1acc-4cdd-97cf-07c84515db63.960381a50bbb9142e9ab989eaff9b9e825477b0b.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b71f158a6df8bc794c22d9aa2ebe3f2247e0ab85/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3cc8a830-1acc-4cdd-97cf-07c84515db63.960381a50bbb9142e9ab989eaff9b9e825477b0b.de-de.xlf", "", "", "3cc8a830-1acc-4cdd-97cf-07c84515db63.960381a50bbb9142e9ab989eaff9b9e825477b0b.de-de.xlf") | Out-Null

$wsDe.Range("H3").Value = "2016-03-19 00:47:40"
